$d = $word.ActiveDocument

# Locate the word "green" inside the "Try changing..." note and replace it
# with "blue" (the surrounding sentence keeps its existing italic formatting).
$rng = $d.Range(0, $d.Content.End)
$found = $rng.Find.Execute("green", $true, $false, $false, $false, $false, `
                            $true, 1, $false, "blue", 2)
if (-not $found) {
    throw "Could not find the word 'green' to replace."
}

# $rng now spans the freshly inserted word "blue". Force a run boundary right
# before it (mirrors Word's own behavior of splitting the run at the edit
# point) using a short-lived bookmark, then drop the helper bookmark again.
$splitPoint = $d.Range($rng.Start, $rng.Start)
$d.Bookmarks.Add("TempSplit", $splitPoint)

# Word re-anchors its hidden "_GoBack" bookmark (collapsed) at the end of the
# most recent edit -- i.e. right after "blue" -- removing it from wherever it
# used to live (the trailing empty paragraph in this document).
$editEnd = $d.Range($rng.End, $rng.End)
$d.Bookmarks.Add("_GoBack", $editEnd)

$d.Bookmarks("TempSplit").Delete()
